$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = 2021
$ws.Range("B32").Value = "大冒険ブック 旅する王国と無名階域"
$ws.Range("C32").Value = "Great Adventure Book: Traveling Kingdom and Anonymous Floors"
$ws.Range("D32").Value = "Kadokawa"
$ws.Range("E32").Value = "great_adventure_book.jpg"
$ws.Range("F32").Value = "supplement"

$ws.Columns.Item(3).ColumnWidth = 53.0

$ws.Range("F33").Select()

